# Applies two changes from the commit:
#  1. The table on slide 6 gets a new (built-in) table style id.
#  2. The deck's theme colour scheme ("Integral") is replaced by the
#     standard "Office Theme" colour scheme (the theme1.xml / theme2.xml
#     parts effectively swap their colour schemes).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 (the "SOURCES OF FINANCE" table) ---------
$slide = $p.Slides.Item(6)
foreach ($shape in $slide.Shapes) {
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{8DCE5809-6F5C-49DB-A0F3-BC83885D1CE8}")
    }
}

# --- 2. Swap the theme colour scheme for the "Office Theme" palette -----
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$colors.Item(1).RGB = 0        # dk1      000000
$colors.Item(2).RGB = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB = 6968388  # dk2      44546A
$colors.Item(4).RGB = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB = 3243501  # accent2  ED7D31
$colors.Item(7).RGB = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB = 49407    # accent4  FFC000
$colors.Item(9).RGB = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456 # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink     0563C1
$colors.Item(12).RGB = 7491477  # folHlink  954F72
